$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3 (shifting old row 3 down to row 4, old row 4 down to row 5)
$ws.Rows.Item(3).Insert()

# New row 3 gets the new data
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 45063
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112035
$ws.Range("G3").Value = "Bruselas (repollito)"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 21000
$ws.Range("L3").Value = 22000
$ws.Range("M3").Value = 21500
$ws.Range("N3").Value = "$/malla 15 kilos"
$ws.Range("O3").Value = "Provincia de Quillota"
$ws.Range("P3").Value = 1433
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = "Hortaliza"
